$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (DÜZENLİ EFT) - İŞBANKASI (E6) was blank, now filled in
$ws.Range("E6").Value = "8.300,01 TL - 199,41 TL"

# Row 12 (GİDEN SWIFT) - FINASNBANK (K12) updated figures
$ws.Range("K12").Value = "WU: ,USD–; Diğer: 529 TL–4.454,74 TL"

# Row 13 (GELEN SWIFT) - İŞBANKASI (E13) and YKB (F13) were blank, now filled in
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 1.114 TL"
$ws.Range("F13").Value = "Hesaba: Asgari 795 TL | Azami 4.005 TL"

# Row 14 (GİDEN SWIFT - Mobil) - İŞBANKASI (E14) and YKB (F14) were blank, now filled in
$ws.Range("E14").Value = "2.170 TL - 2.170 TL"
$ws.Range("F14").Value = "2.785,72 TL - 12.380,95 TL"

# Row 14 (GİDEN SWIFT - Mobil) - FINASNBANK (K14) updated figures
$ws.Range("K14").Value = "1.196,51 TL - 5.583,74 TL"
